$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2,8).Value = 282.66666
$ws.Cells.Item(2,9).Value = 279.3125
$ws.Cells.Item(2,10).Value = 309.5
$ws.Cells.Item(2,11).Value = 279.3125
$ws.Cells.Item(2,12).Value = 309.5
$ws.Cells.Item(2,13).Value = -166.3125
$ws.Cells.Item(2,14).Value = -535.5
$ws.Cells.Item(18,8).Value = 9078.571
$ws.Cells.Item(18,9).Value = 10275
$ws.Cells.Item(18,10).Value = 1900
$ws.Cells.Item(18,11).Value = 10275
$ws.Cells.Item(18,12).Value = 1900
$ws.Cells.Item(18,13).Value = -9991
$ws.Cells.Item(18,14).Value = -2468
$ws.Cells.Item(70,8).Value = 3528.875
$ws.Cells.Item(70,10).Value = 3855.3333
$ws.Cells.Item(70,12).Value = 11565.9999
$ws.Cells.Item(70,14).Value = -12105.9999
$ws.Cells.Item(73,8).Value = 3528.875
$ws.Cells.Item(73,10).Value = 3855.3333
$ws.Cells.Item(73,12).Value = 11565.9999
$ws.Cells.Item(73,14).Value = -13437.9999
$ws.Cells.Item(94,8).Value = 1108
$ws.Cells.Item(94,9).Value = 929.7
$ws.Cells.Item(94,11).Value = 929.7
$ws.Cells.Item(94,13).Value = -478.7
$ws.Cells.Item(103,8).Value = 1542.4286
$ws.Cells.Item(103,9).Value = 966.3333
$ws.Cells.Item(103,10).Value = 1974.5
$ws.Cells.Item(103,11).Value = 2898.9999
$ws.Cells.Item(103,12).Value = 5923.5
$ws.Cells.Item(103,13).Value = -2312.9999
$ws.Cells.Item(103,14).Value = -7095.5
$ws.Cells.Item(112,8).Value = 1526.4773
$ws.Cells.Item(112,10).Value = 1597.9
$ws.Cells.Item(112,12).Value = 4793.700000000001
$ws.Cells.Item(112,14).Value = -7009.700000000001
$ws.Cells.Item(138,8).Value = 2426.2083
$ws.Cells.Item(138,9).Value = 1366.5714
$ws.Cells.Item(138,10).Value = 3909.7
$ws.Cells.Item(138,11).Value = 4099.7142
$ws.Cells.Item(138,12).Value = 11729.1
$ws.Cells.Item(138,13).Value = 1040.2858
$ws.Cells.Item(138,14).Value = -22009.1
$ws.Cells.Item(141,8).Value = 2839.8333
$ws.Cells.Item(141,9).Value = 2839.8333
$ws.Cells.Item(141,11).Value = 8519.499899999999
$ws.Cells.Item(141,13).Value = -3339.499899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74,8).Value = 2033.5238
$ws.Cells.Item(74,9).Value = 2061.611
$ws.Cells.Item(74,11).Value = 2061.611
$ws.Cells.Item(74,13).Value = -1187.611
$ws.Cells.Item(77,8).Value = 2033.5238
$ws.Cells.Item(77,9).Value = 2061.611
$ws.Cells.Item(77,11).Value = 10308.055
$ws.Cells.Item(77,13).Value = -5940.055
$ws.Cells.Item(130,8).Value = 189283.3
$ws.Cells.Item(130,10).Value = 189283.3
$ws.Cells.Item(130,12).Value = 189283.3
$ws.Cells.Item(130,14).Value = -199323.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81,8).Value = 49273.43
$ws.Cells.Item(81,10).Value = 54986.5
$ws.Cells.Item(81,12).Value = 54986.5
$ws.Cells.Item(81,14).Value = -57108.5
$ws.Cells.Item(84,8).Value = 49273.43
$ws.Cells.Item(84,10).Value = 54986.5
$ws.Cells.Item(84,12).Value = 164959.5
$ws.Cells.Item(84,14).Value = -175567.5
$ws.Cells.Item(138,8).Value = 72164.664
$ws.Cells.Item(138,10).Value = 72234.14999999999
$ws.Cells.Item(138,12).Value = 72234.14999999999
$ws.Cells.Item(138,14).Value = -82514.14999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22,8).Value = 1100
$ws.Cells.Item(22,9).Value = 200
$ws.Cells.Item(22,11).Value = 200
$ws.Cells.Item(22,13).Value = 150
$ws.Cells.Item(31,8).Value = 3077.4866
$ws.Cells.Item(31,9).Value = 1989.6666
$ws.Cells.Item(31,11).Value = 1989.6666
$ws.Cells.Item(31,13).Value = -1694.6666
$ws.Cells.Item(34,8).Value = 3077.4866
$ws.Cells.Item(34,9).Value = 1989.6666
$ws.Cells.Item(34,11).Value = 1989.6666
$ws.Cells.Item(34,13).Value = -1787.6666
$ws.Cells.Item(41,8).Value = 293.66666
$ws.Cells.Item(41,9).Value = 293.66666
$ws.Cells.Item(41,11).Value = 293.66666
$ws.Cells.Item(41,13).Value = 134.33334
$ws.Cells.Item(62,8).Value = 3684.6365
$ws.Cells.Item(62,9).Value = 4104.5557
$ws.Cells.Item(62,10).Value = 1795
$ws.Cells.Item(62,11).Value = 4104.5557
$ws.Cells.Item(62,12).Value = 1795
$ws.Cells.Item(62,13).Value = -3480.5557
$ws.Cells.Item(62,14).Value = -3043
$ws.Cells.Item(65,8).Value = 3684.6365
$ws.Cells.Item(65,9).Value = 4104.5557
$ws.Cells.Item(65,10).Value = 1795
$ws.Cells.Item(65,11).Value = 20522.7785
$ws.Cells.Item(65,12).Value = 8975
$ws.Cells.Item(65,13).Value = -17402.7785
$ws.Cells.Item(65,14).Value = -15215
$ws.Cells.Item(86,8).Value = 18346.5
$ws.Cells.Item(86,9).Value = 23693.2
$ws.Cells.Item(86,11).Value = 23693.2
$ws.Cells.Item(86,13).Value = -22570.2
$ws.Cells.Item(89,8).Value = 18346.5
$ws.Cells.Item(89,9).Value = 23693.2
$ws.Cells.Item(89,11).Value = 118466
$ws.Cells.Item(89,13).Value = -112850
$ws.Cells.Item(99,8).Value = 8636250
$ws.Cells.Item(99,10).Value = 18189846
$ws.Cells.Item(99,12).Value = 18189846
$ws.Cells.Item(99,14).Value = -18192842
$ws.Cells.Item(126,8).Value = 8636250
$ws.Cells.Item(126,10).Value = 18189846
$ws.Cells.Item(126,12).Value = 54569538
$ws.Cells.Item(126,14).Value = -54574478
$ws.Cells.Item(132,8).Value = 2087.1072
$ws.Cells.Item(132,9).Value = 1377.6
$ws.Cells.Item(132,11).Value = 4132.799999999999
$ws.Cells.Item(132,13).Value = -1602.799999999999
$ws.Cells.Item(134,8).Value = 3900.1785
$ws.Cells.Item(134,9).Value = 2318.342
$ws.Cells.Item(134,11).Value = 6955.026
$ws.Cells.Item(134,13).Value = -4420.026
$ws.Cells.Item(135,8).Value = 67995.96000000001
$ws.Cells.Item(135,10).Value = 67995.96000000001
$ws.Cells.Item(135,12).Value = 67995.96000000001
$ws.Cells.Item(135,14).Value = -78135.96000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131,8).Value = 4448.773
$ws.Cells.Item(131,9).Value = 3218.2856
$ws.Cells.Item(131,11).Value = 9654.856800000001
$ws.Cells.Item(131,13).Value = -4614.856800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113,8).Value = 6122
$ws.Cells.Item(113,10).Value = 8837.799999999999
$ws.Cells.Item(113,12).Value = 8837.799999999999
$ws.Cells.Item(113,14).Value = -13177.8
$ws.Cells.Item(122,8).Value = 5034.727
$ws.Cells.Item(122,9).Value = 3340.2856
$ws.Cells.Item(122,11).Value = 10020.8568
$ws.Cells.Item(122,13).Value = -7570.856800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(25,8).Value = 1000
$ws.Cells.Item(25,9).Value = 1000
$ws.Cells.Item(25,10).Value = 0
$ws.Cells.Item(25,11).Value = 1000
$ws.Cells.Item(25,12).Value = 0
$ws.Cells.Item(25,13).ClearContents()
$ws.Cells.Item(25,14).Value = -770
$ws.Cells.Item(40,8).Value = 12554.467
$ws.Cells.Item(40,9).Value = 51594
$ws.Cells.Item(40,11).Value = 51594
$ws.Cells.Item(40,13).Value = -51458
$ws.Cells.Item(55,8).Value = 2186
$ws.Cells.Item(55,9).Value = 3591.6667
$ws.Cells.Item(55,11).Value = 3591.6667
$ws.Cells.Item(55,13).Value = -3418.6667
$ws.Cells.Item(68,8).Value = 266143.75
$ws.Cells.Item(68,9).Value = 169053.17
$ws.Cells.Item(68,11).Value = 169053.17
$ws.Cells.Item(68,13).Value = -168304.17
$ws.Cells.Item(71,8).Value = 266143.75
$ws.Cells.Item(71,9).Value = 169053.17
$ws.Cells.Item(71,11).Value = 845265.8500000001
$ws.Cells.Item(71,13).Value = -841521.8500000001
$ws.Cells.Item(82,8).Value = 4124.2104
$ws.Cells.Item(82,9).Value = 2298.3333
$ws.Cells.Item(82,10).Value = 7254.2856
$ws.Cells.Item(82,11).Value = 2298.3333
$ws.Cells.Item(82,12).Value = 7254.2856
$ws.Cells.Item(82,13).Value = -1937.3333
$ws.Cells.Item(82,14).Value = -7976.2856
$ws.Cells.Item(85,8).Value = 4124.2104
$ws.Cells.Item(85,9).Value = 2298.3333
$ws.Cells.Item(85,10).Value = 7254.2856
$ws.Cells.Item(85,11).Value = 2298.3333
$ws.Cells.Item(85,12).Value = 7254.2856
$ws.Cells.Item(85,13).Value = -1050.3333
$ws.Cells.Item(85,14).Value = -9750.285599999999
$ws.Cells.Item(132,8).Value = 3777.8088
$ws.Cells.Item(132,9).Value = 2693.638
$ws.Cells.Item(132,11).Value = 8080.914
$ws.Cells.Item(132,13).Value = -5550.914
$ws.Cells.Item(136,8).Value = 4730.5356
$ws.Cells.Item(136,9).Value = 2598.6428
$ws.Cells.Item(136,11).Value = 7795.928400000001
$ws.Cells.Item(136,13).Value = -5245.928400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132,8).Value = 2070.5862
$ws.Cells.Item(132,9).Value = 1022.4211
$ws.Cells.Item(132,11).Value = 2070.5862
$ws.Cells.Item(132,13).Value = -537.2633000000001
